$d = $word.ActiveDocument
$nl = [char]11

function Replace-ParagraphText($para, [string]$newText) {
    $r = $para.Range
    $old = $r.Text
    if ($old.Length -gt 0 -and [int][char]$old[$old.Length-1] -eq 13) {
        $old = $old.Substring(0, $old.Length - 1)
    }
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Replace failed for paragraph text: $old"
    }
}

# --- "Why Now?" bullet list: Key factors include ---
$newText10 = "Key factors include:" + $nl + "• Trading card and collectibles markets are expanding, with strong retail demand and investment interest [1][6]." + $nl + "• Real-time marketplace data and APIs enable live pricing aggregation." + $nl + "• Counterfeits in secondary markets increase the need for authenticity verification (supports an AI-first approach)." + $nl + "• AI maturity: Computer vision (Rekognition) + LLM reasoning (Bedrock) enable explainable evaluation workflows." + $nl + "• Pokémon TCG content and digital engagement remain high, sustaining volatility and attention [7]."
Replace-ParagraphText $d.Paragraphs.Item(10) $newText10

$newText13 = "• Retail signals: Walmart reported trading card sales up ~200% (Feb 2024 → Jun 2025) and Pokémon card sales >10× YoY; broader retail momentum persists [6]."
Replace-ParagraphText $d.Paragraphs.Item(13) $newText13

$newText59 = "• Gamification and digital twin collectibles (e.g., Pokémon TCG Pocket surpassed 100M downloads in Feb 2025) [7]."
Replace-ParagraphText $d.Paragraphs.Item(59) $newText59

$newText61 = "Reddit:" + $nl + "• r/PokemonTCG (large, active community for pulls, decks, and pricing) [8]." + $nl + "• r/PokeInvesting (active market & investment discussion) [9]."
Replace-ParagraphText $d.Paragraphs.Item(61) $newText61

# --- References section (new content appended at end of document) ---
$headingPara = $d.Paragraphs.Add()
$headingPara.Style = "Heading1"
$headingPara.Range.Text = "References"

function Add-ReferenceParagraph([string]$citation, [string]$url) {
    $para = $d.Paragraphs.Add()
    # Paragraphs.Add() inherits the style of the preceding paragraph (e.g. Heading1);
    # force it back to Normal for these body/reference paragraphs.
    $para.Style = "Normal"
    $para.Range.Text = $citation
    # Toggle bold True->False so the run explicitly records b=0 (matches source authoring tool)
    $para.Range.Font.Bold = $true
    $para.Range.Font.Bold = $false

    $insPos = $para.Range.End - 1
    $insRange = $d.Range($insPos, $insPos)
    $insRange.InsertAfter($url)
    $urlRange = $d.Range($insPos, $insPos + $url.Length)
    # wdUndefined clears the explicit Bold flag so the 2nd run has no <w:b> override
    $urlRange.Font.Bold = 9999999
}

Add-ReferenceParagraph "[1] Grand View Research – Collectibles Market size `$294.23B (2023) to `$422.56B (2030), CAGR 5.5%.  " "https://www.grandviewresearch.com/industry-analysis/collectibles-market-report"
Add-ReferenceParagraph "[2] Zion Market Research – Trading Card Game Market `$7.43B (2024) to `$15.84B (2034), CAGR 7.86%.  " "https://www.zionmarketresearch.com/report/trading-card-game-market"
Add-ReferenceParagraph "[3] Verified Market Research – Sports Trading Card Market `$12.62B (2024) to `$23.08B (2031), CAGR 7.8%.  " "https://www.verifiedmarketresearch.com/product/sports-trading-card-market/"
Add-ReferenceParagraph "[4] Astute Analytica via Yahoo Finance – TCG Authentication Services `$2.24B (2024) to `$6.61B (2033), CAGR 13.1%.  " "https://finance.yahoo.com/news/global-trading-card-game-authentication-163000942.html"
Add-ReferenceParagraph "[5] The Economic Times (citing Washington Post & Card Ladder) – Pokémon cards ~3,821% cumulative return since 2004; S&P 500 ~483%.  " "https://m.economictimes.com/news/international/global-trends/pikachus-wild-run-pokemon-cards-give-3821-return-to-blow-past-sp-500s-483/articleshow/123823347.cms"
Add-ReferenceParagraph "[6] Axios – Retail momentum: trading card sales up at Walmart (200% overall; Pokémon >10x YoY).  " "https://www.axios.com/2025/08/23/pokemon-card-trading-cards-boom-target-ebay-sales"
Add-ReferenceParagraph "[7] PocketGamer.biz – Pokémon TCG Pocket surpasses 100M downloads (Feb 28, 2025).  " "https://www.pocketgamer.biz/pokmon-tcg-pocket-surpasses-100-million-downloads/"
Add-ReferenceParagraph "[8] r/PokemonTCG subreddit (active community).  " "https://www.reddit.com/r/PokemonTCG/"
Add-ReferenceParagraph "[9] r/PokeInvesting subreddit (active community).  " "https://www.reddit.com/r/PokeInvesting/"
